# Sync up to the latest carina archetype
# - Refresh the demo phone records on the "GSMArena" sheet (rows 2-4,
#   columns D:H -> model/display/camera/ram/battery) with the new sample
#   data (Galaxy S10+, Galaxy Fold, Galaxy M10).
# - Leave the "Calculator" sheet's data untouched.
# - The active sheet/selection moves to GSMArena!H6 (was Calculator!E4).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GSMArena")

# Row 2: Samsung Galaxy S10+
$ws1.Range("D2").Value = "Galaxy S10+"
$ws1.Range("E2").Value = "6.4"""
$ws1.Range("F2").Value = "16MP"
$ws1.Range("G2").Value = "12GB RAM"
$ws1.Range("H2").Value = "4100mAh"

# Row 3: Samsung Galaxy Fold
$ws1.Range("D3").Value = "Galaxy Fold"
$ws1.Range("E3").Value = "7.3"""
$ws1.Range("F3").Value = "16MP"
$ws1.Range("G3").Value = "12GB RAM"
$ws1.Range("H3").Value = "4380mAh"

# Row 4: Samsung Galaxy M10
$ws1.Range("D4").Value = "Galaxy M10"
$ws1.Range("E4").Value = "6.22"""
$ws1.Range("F4").Value = "13MP"
$ws1.Range("G4").Value = "3GB RAM"
$ws1.Range("H4").Value = "3400mAh"

# GSMArena becomes the active/selected sheet with H6 selected.
$ws1.Activate()
$ws1.Range("H6").Select()
